$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet from "Sheet1" to "Table"
$ws.Name = "Table"

# Select the sheet and update the active cell / selection to E12
$ws.Activate()
$ws.Range("E12").Select()
